$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name and Link for rows 18-19 (Litecoin/Uniswap swapped)
$ws.Range("B18").Value = "Uniswap"
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"

# Update Price column (D) - force text to preserve formatting (leading apostrophe where the
# value would otherwise be auto-converted to a number by Excel)
$ws.Range("D2").Value = "37.869.15"
$ws.Range("D3").Value = "2.083.87"
$ws.Range("D5").Value = "'233.13"
$ws.Range("D6").Value = "'0.625"
$ws.Range("D7").Value = "'59.09"
$ws.Range("D9").Value = "'0.394"
$ws.Range("D10").Value = "'0.0786"
$ws.Range("D12").Value = "'14.72"
$ws.Range("D13").Value = "'21.19"
$ws.Range("D14").Value = "'0.773"
$ws.Range("D15").Value = "'5.33"
$ws.Range("D16").Value = "2.059.94"
$ws.Range("D17").Value = "37.800.67"
$ws.Range("D18").Value = "'6.12"
$ws.Range("D19").Value = "'71.57"
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("D21").Value = "'228.03"
$ws.Range("D25").Value = "'171.81"
$ws.Range("D26").Value = "'9.21"
$ws.Range("D32").Value = "'4.71"
$ws.Range("D33").Value = "'0.0630"
$ws.Range("D38").Value = "'5.41"
$ws.Range("D39").Value = "'0.0981"
$ws.Range("D40").Value = "'99.08"
$ws.Range("D44").Value = "1.445.68"
$ws.Range("D50").Value = "2.277.31"
$ws.Range("D51").Value = "'46.82"

# Update Volume(1h) column (E)
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  +8.05%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +0.68%  "
